# Update countries & provincias Spain
# - refresh "last updated" timestamp
# - refresh case counters for the countries whose stats moved
# - a handful of countries swapped ranking position (same row stays the
#   same rank slot, but which country occupies it changed), so those rows
#   get both a new country name (col A) and new stats (cols B-H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 8 de Septiembre de 2020 a las 02:32'

$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 6485567
$ws.Cells.Item(4, 3).Value = 25317
$ws.Cells.Item(4, 4).Value = 3758618
$ws.Cells.Item(4, 5).Value = 2533415
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 286
$ws.Cells.Item(4, 8).Value = 193534

$ws.Cells.Item(6, 1).Value = 'Brasil'
$ws.Cells.Item(6, 2).Value = 4147794
$ws.Cells.Item(6, 3).Value = 10188
$ws.Cells.Item(6, 4).Value = 3355564
$ws.Cells.Item(6, 5).Value = 665229
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 315
$ws.Cells.Item(6, 8).Value = 127001

$ws.Cells.Item(36, 1).Value = 'Panama'
$ws.Cells.Item(36, 2).Value = 97578
$ws.Cells.Item(36, 3).Value = 535
$ws.Cells.Item(36, 4).Value = 70247
$ws.Cells.Item(36, 5).Value = 25232
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 13
$ws.Cells.Item(36, 8).Value = 2099

$ws.Cells.Item(56, 1).Value = 'Venezuela'
$ws.Cells.Item(56, 2).Value = 54350
$ws.Cells.Item(56, 3).Value = 1061
$ws.Cells.Item(56, 4).Value = 43753
$ws.Cells.Item(56, 5).Value = 10161
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 8
$ws.Cells.Item(56, 8).Value = 436

$ws.Cells.Item(72, 1).Value = 'Chequia'
$ws.Cells.Item(72, 2).Value = 28716
$ws.Cells.Item(72, 3).Value = 560
$ws.Cells.Item(72, 4).Value = 19855
$ws.Cells.Item(72, 5).Value = 8424
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 437

$ws.Cells.Item(111, 1).Value = 'Tunez'
$ws.Cells.Item(111, 2).Value = 5124
$ws.Cells.Item(111, 3).Value = 83
$ws.Cells.Item(111, 4).Value = 1788
$ws.Cells.Item(111, 5).Value = 3242
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 94

$ws.Cells.Item(113, 1).Value = 'Congo'
$ws.Cells.Item(113, 2).Value = 4891
$ws.Cells.Item(113, 3).Value = 263
$ws.Cells.Item(113, 4).Value = 3887
$ws.Cells.Item(113, 5).Value = 902
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 102

$ws.Cells.Item(114, 1).Value = 'Hong Kong'
$ws.Cells.Item(114, 2).Value = 4890
$ws.Cells.Item(114, 3).Value = 11
$ws.Cells.Item(114, 4).Value = 4524
$ws.Cells.Item(114, 5).Value = 268
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 4
$ws.Cells.Item(114, 8).Value = 98

$ws.Cells.Item(115, 1).Value = 'Suazilandia'
$ws.Cells.Item(115, 2).Value = 4884
$ws.Cells.Item(115, 3).Value = 31
$ws.Cells.Item(115, 4).Value = 4029
$ws.Cells.Item(115, 5).Value = 761
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 94

$ws.Cells.Item(116, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(116, 2).Value = 4729
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 1818
$ws.Cells.Item(116, 5).Value = 2849
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 62

$ws.Cells.Item(117, 1).Value = 'Nicaragua'
$ws.Cells.Item(117, 2).Value = 4668
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = 2913
$ws.Cells.Item(117, 5).Value = 1614
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 141

$ws.Cells.Item(118, 1).Value = 'Eslovaquia'
$ws.Cells.Item(118, 2).Value = 4636
$ws.Cells.Item(118, 3).Value = 22
$ws.Cells.Item(118, 4).Value = 2836
$ws.Cells.Item(118, 5).Value = 1763
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 37

$ws.Cells.Item(121, 1).Value = 'Surinam'
$ws.Cells.Item(121, 2).Value = 4360
$ws.Cells.Item(121, 3).Value = 14
$ws.Cells.Item(121, 4).Value = 3544
$ws.Cells.Item(121, 5).Value = 725
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 6
$ws.Cells.Item(121, 8).Value = 91

$ws.Cells.Item(122, 1).Value = 'Cabo Verde'
$ws.Cells.Item(122, 2).Value = 4358
$ws.Cells.Item(122, 3).Value = 28
$ws.Cells.Item(122, 4).Value = 3790
$ws.Cells.Item(122, 5).Value = 526
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 42

$ws.Cells.Item(123, 1).Value = 'Cuba'
$ws.Cells.Item(123, 2).Value = 4352
$ws.Cells.Item(123, 3).Value = 43
$ws.Cells.Item(123, 4).Value = 3642
$ws.Cells.Item(123, 5).Value = 608
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 102

$ws.Cells.Item(128, 1).Value = 'Siria'
$ws.Cells.Item(128, 2).Value = 3229
$ws.Cells.Item(128, 3).Value = 58
$ws.Cells.Item(128, 4).Value = 744
$ws.Cells.Item(128, 5).Value = 2348
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 3
$ws.Cells.Item(128, 8).Value = 137

$ws.Cells.Item(129, 1).Value = 'Gambia'
$ws.Cells.Item(129, 2).Value = 3197
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 1315
$ws.Cells.Item(129, 5).Value = 1783
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 99

$ws.Cells.Item(130, 1).Value = 'Eslovenia'
$ws.Cells.Item(130, 2).Value = 3190
$ws.Cells.Item(130, 3).Value = 25
$ws.Cells.Item(130, 4).Value = 2530
$ws.Cells.Item(130, 5).Value = 525
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 135

$ws.Cells.Item(131, 1).Value = 'Sri Lanka'
$ws.Cells.Item(131, 2).Value = 3123
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 2926
$ws.Cells.Item(131, 5).Value = 185
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 12

$ws.Cells.Item(145, 1).Value = 'Islandia'
$ws.Cells.Item(145, 2).Value = 2143
$ws.Cells.Item(145, 3).Value = 2
$ws.Cells.Item(145, 4).Value = 2057
$ws.Cells.Item(145, 5).Value = 76
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 10

$ws.Cells.Item(148, 1).Value = 'Sierra Leona'
$ws.Cells.Item(148, 2).Value = 2055
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(148, 4).Value = 1611
$ws.Cells.Item(148, 5).Value = 373
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 71

$ws.Cells.Item(153, 1).Value = 'Guyana'
$ws.Cells.Item(153, 2).Value = 1560
$ws.Cells.Item(153, 3).Value = 92
$ws.Cells.Item(153, 4).Value = 962
$ws.Cells.Item(153, 5).Value = 551
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 1
$ws.Cells.Item(153, 8).Value = 47

$ws.Cells.Item(154, 1).Value = 'Birmania'
$ws.Cells.Item(154, 2).Value = 1518
$ws.Cells.Item(154, 3).Value = 99
$ws.Cells.Item(154, 4).Value = 388
$ws.Cells.Item(154, 5).Value = 1122
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 8

$ws.Cells.Item(155, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(155, 2).Value = 1510
$ws.Cells.Item(155, 3).Value = 1
$ws.Cells.Item(155, 4).Value = 1237
$ws.Cells.Item(155, 5).Value = 251
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 1
$ws.Cells.Item(155, 8).Value = 22

$ws.Cells.Item(156, 1).Value = 'Togo'
$ws.Cells.Item(156, 2).Value = 1493
$ws.Cells.Item(156, 3).Value = 5
$ws.Cells.Item(156, 4).Value = 1114
$ws.Cells.Item(156, 5).Value = 346
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 1
$ws.Cells.Item(156, 8).Value = 33

$ws.Cells.Item(157, 1).Value = 'Burkina Faso'
$ws.Cells.Item(157, 2).Value = 1463
$ws.Cells.Item(157, 3).Value = 11
$ws.Cells.Item(157, 4).Value = 1112
$ws.Cells.Item(157, 5).Value = 295
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 1
$ws.Cells.Item(157, 8).Value = 56

$ws.Cells.Item(161, 1).Value = 'Belice'
$ws.Cells.Item(161, 2).Value = 1307
$ws.Cells.Item(161, 3).Value = 46
$ws.Cells.Item(161, 4).Value = 314
$ws.Cells.Item(161, 5).Value = 977
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 16

$ws.Cells.Item(162, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(162, 2).Value = 1261
$ws.Cells.Item(162, 3).Value = 46
$ws.Cells.Item(162, 4).Value = 934
$ws.Cells.Item(162, 5).Value = 274
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 53
